$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 19:38"

$ws.Range("B4").Value = 4936836
$ws.Range("C4").Value = 18416
$ws.Range("D4").Value = 2500773
$ws.Range("E4").Value = 2275230
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 543
$ws.Range("H4").Value = 160833

$ws.Range("B6").Value = 1959822
$ws.Range("C6").Value = 53209
$ws.Range("D6").Value = 1325040
$ws.Range("E6").Value = 594050
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 912
$ws.Range("H6").Value = 40732

$ws.Range("B20").Value = 236112
$ws.Range("C20").Value = 1178
$ws.Range("D20").Value = 219506
$ws.Range("E20").Value = 10822
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 5784

$ws.Range("A21").Value = "Alemania"
$ws.Range("B21").Value = 213613
$ws.Range("C21").Value = 533
$ws.Range("D21").Value = 194700
$ws.Range("E21").Value = 9671
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 10
$ws.Range("H21").Value = 9242

$ws.Range("A22").Value = "Argentina"
$ws.Range("B22").Value = 213535
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 96948
$ws.Range("E22").Value = 112578
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 30
$ws.Range("H22").Value = 4009

$ws.Range("B25").Value = 118037
$ws.Range("C25").Value = 245
$ws.Range("D25").Value = 102599
$ws.Range("E25").Value = 6478
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 8960

$ws.Range("B31").Value = 88866
$ws.Range("C31").Value = 903
$ws.Range("D31").Value = 71168
$ws.Range("E31").Value = 11851
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 39
$ws.Range("H31").Value = 5847

$ws.Range("B36").Value = 77595
$ws.Range("C36").Value = 1397
$ws.Range("D36").Value = 51378
$ws.Range("E36").Value = 25652
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 4
$ws.Range("H36").Value = 565

$ws.Range("A61").Value = "Marruecos"
$ws.Range("B61").Value = 28500
$ws.Range("C61").Value = 1283
$ws.Range("D61").Value = 19994
$ws.Range("E61").Value = 8071
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 18
$ws.Range("H61").Value = 435

$ws.Range("A62").Value = "Uzbekistan"
$ws.Range("B62").Value = 27554
$ws.Range("C62").Value = 507
$ws.Range("D62").Value = 18682
$ws.Range("E62").Value = 8703
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 169

$ws.Range("B64").Value = 26303
$ws.Range("C64").Value = 50
$ws.Range("D64").Value = 23364
$ws.Range("E64").Value = 1176
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 1763

$ws.Range("B89").Value = 8069
$ws.Range("C89").Value = 71
$ws.Range("D89").Value = 7075
$ws.Range("E89").Value = 947
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 47

$ws.Range("B96").Value = 7007
$ws.Range("C96").Value = 90
$ws.Range("D96").Value = 5623
$ws.Range("E96").Value = 1266
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 118

$ws.Range("A106").Value = "Maldivas"
$ws.Range("B106").Value = 4594
$ws.Range("C106").Value = 148
$ws.Range("D106").Value = 2703
$ws.Range("E106").Value = 1872
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 19

$ws.Range("A107").Value = "Hungria"
$ws.Range("B107").Value = 4564
$ws.Range("C107").Value = 11
$ws.Range("D107").Value = 3431
$ws.Range("E107").Value = 534
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 599

$ws.Range("B130").Value = 2079
$ws.Range("C130").Value = 50
$ws.Range("D130").Value = 778
$ws.Range("E130").Value = 1286
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 15

$ws.Range("B135").Value = 1860
$ws.Range("C135").Value = 5
$ws.Range("D135").Value = 1401
$ws.Range("E135").Value = 392
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 67

$ws.Range("B142").Value = 1231
$ws.Range("C142").Value = 7
$ws.Range("D142").Value = 1160
$ws.Range("E142").Value = 60
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 11

$ws.Range("A181").Value = "Aruba"
$ws.Range("B181").Value = 171
$ws.Range("C181").Value = 39
$ws.Range("D181").Value = 112
$ws.Range("E181").Value = 56
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 3

$ws.Range("A182").Value = "Bermudas"
$ws.Range("B182").Value = 157
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 144
$ws.Range("E182").Value = 4
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 9

$ws.Range("A183").Value = "Papua Nueva Guinea"
$ws.Range("B183").Value = 153
$ws.Range("C183").Value = 39
$ws.Range("D183").Value = 44
$ws.Range("E183").Value = 107
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 2

$ws.Range("A184").Value = "San Martin (Parte Holandesa)"
$ws.Range("B184").Value = 150
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 64
$ws.Range("E184").Value = 70
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 16

$ws.Range("A185").Value = "Brunei"
$ws.Range("B185").Value = 141
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 138
$ws.Range("E185").Value = 0
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 3

$ws.Range("A186").Value = "Barbados"
$ws.Range("B186").Value = 132
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 98
$ws.Range("E186").Value = 27
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 7
